# SE-CPE4A-19-20-ATTENDANCE.xlsx edit
# Fills in attendance records for the first week of classes that were
# still left blank, adds two more class dates to the header row, and
# leaves the selection on the cell the author was last working on.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 holds the class-session dates; two more sessions (Jan 29 & 30, 2020)
# are recorded, replacing the blank placeholders in G6/H6.
$ws.Range("G6").Value = 43859
$ws.Range("H6").Value = 43860

# Attendance marks (1 = present, 0 = absent) for columns E-H, rows 7-54.
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 1
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 1
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 1
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = 1
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = 0
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 1
$ws.Range("H21").Value = 1
$ws.Range("E22").Value = 1
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 1
$ws.Range("E23").Value = 1
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 1
$ws.Range("E24").Value = 1
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = 1
$ws.Range("H24").Value = 1
$ws.Range("E25").Value = 1
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 1
$ws.Range("H25").Value = 1
$ws.Range("E26").Value = 1
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 0
$ws.Range("E28").Value = 1
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 1
$ws.Range("E29").Value = 1
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 1
$ws.Range("H29").Value = 1
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = 1
$ws.Range("H30").Value = 1
$ws.Range("E31").Value = 1
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 1
$ws.Range("E32").Value = 1
$ws.Range("F32").Value = 1
$ws.Range("G32").Value = 1
$ws.Range("H32").Value = 1
$ws.Range("E33").Value = 1
$ws.Range("F33").Value = 1
$ws.Range("G33").Value = 1
$ws.Range("H33").Value = 1
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = 1
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 1
$ws.Range("E35").Value = 1
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = 1
$ws.Range("E36").Value = 1
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 1
$ws.Range("H36").Value = 1
$ws.Range("E37").Value = 1
$ws.Range("F37").Value = 1
$ws.Range("G37").Value = 1
$ws.Range("H37").Value = 1
$ws.Range("E38").Value = 1
$ws.Range("F38").Value = 1
$ws.Range("G38").Value = 1
$ws.Range("H38").Value = 1
$ws.Range("E39").Value = 0
$ws.Range("F39").Value = 1
$ws.Range("G39").Value = 1
$ws.Range("H39").Value = 0
$ws.Range("E40").Value = 1
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 1
$ws.Range("H40").Value = 0
$ws.Range("E41").Value = 1
$ws.Range("F41").Value = 1
$ws.Range("G41").Value = 1
$ws.Range("H41").Value = 1
$ws.Range("E42").Value = 1
$ws.Range("F42").Value = 1
$ws.Range("G42").Value = 1
$ws.Range("H42").Value = 1
$ws.Range("E43").Value = 1
$ws.Range("F43").Value = 1
$ws.Range("G43").Value = 1
$ws.Range("H43").Value = 1
$ws.Range("E44").Value = 0
$ws.Range("F44").Value = 1
$ws.Range("G44").Value = 1
$ws.Range("H44").Value = 0
$ws.Range("E45").Value = 1
$ws.Range("F45").Value = 1
$ws.Range("G45").Value = 1
$ws.Range("H45").Value = 1
$ws.Range("E46").Value = 1
$ws.Range("F46").Value = 1
$ws.Range("G46").Value = 1
$ws.Range("H46").Value = 1
$ws.Range("E47").Value = 1
$ws.Range("F47").Value = 1
$ws.Range("G47").Value = 1
$ws.Range("H47").Value = 1
$ws.Range("E48").Value = 1
$ws.Range("F48").Value = 1
$ws.Range("G48").Value = 1
$ws.Range("H48").Value = 1
$ws.Range("E49").Value = 1
$ws.Range("F49").Value = 1
$ws.Range("G49").Value = 1
$ws.Range("H49").Value = 1
$ws.Range("E50").Value = 1
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 1
$ws.Range("H50").Value = 1
$ws.Range("E51").Value = 1
$ws.Range("F51").Value = 1
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 1
$ws.Range("E52").Value = 1
$ws.Range("F52").Value = 1
$ws.Range("G52").Value = 1
$ws.Range("H52").Value = 1
$ws.Range("E53").Value = 1
$ws.Range("F53").Value = 1
$ws.Range("G53").Value = 1
$ws.Range("H53").Value = 1
$ws.Range("E54").Value = 1
$ws.Range("F54").Value = 1
$ws.Range("G54").Value = 1
$ws.Range("H54").Value = 1

# Leave the active selection where the author left it.
$ws.Range("H34").Select()
